# Gantt.xlsx - "Card animation" update
# Updates task assignees (column C) and progress percentages (column D)
# on the "Planificateur Gantt" sheet, plus view/zoom/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Progress (column D) updates ---------------------------------------
$ws.Range("D19").Value = 1      # Base de l'application
$ws.Range("D20").Value = 1      # Pages vitrines de l'application
$ws.Range("D21").Value = 0.8    # Pages sécurisées pour les membres
$ws.Range("D22").Value = 0.6    # Design CSS

$ws.Range("C24").Value = "Vincent, Théo"
$ws.Range("D24").Value = 1      # Base de données

$ws.Range("C25").Value = "Nicolas, Vincent, Théo"
$ws.Range("D25").Value = 0.6    # CRUD

$ws.Range("C26").Value = "Nicolas, Vincent, Théo"
$ws.Range("D26").Value = 0.8    # Sécurisation des accès

$ws.Range("C27").Value = "Nicolas, Vincent, Théo"
$ws.Range("D27").Value = 0.6    # Intéractions des éléments

$ws.Range("C33").Value = "Théo"
$ws.Range("D33").Value = 0.2    # Tests Unitaires

# --- View state: zoom back to 100%, scroll down, reselect -------------
$excel.ActiveWindow.Zoom = 100
$ws.Range("A23").Select()
$ws.Range("D34").Select()
